# Refresh crypto price/volume data (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain-text values such as
# '26.879.66', '0.9971' or '  +1.72%  '. Each refreshed cell is marked as
# Text before the write so Excel doesn't reinterpret the numeric-looking
# strings as numbers (which would drop significant trailing zeros, etc.).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.879.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.731.57'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9971'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.30'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4912'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2603'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06228'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.735.38'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.09'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06910'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6116'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9982'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.640.89'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9973'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007187'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.959.19'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.445'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.577'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.138'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.74'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.33'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.792'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.13%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.18'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.948'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07997'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.684'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04531'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9970'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.604'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.009'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6250'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9359'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.055'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.445'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01505'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.645'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.80'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3874'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.951'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05387'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.925'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.29'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.246'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.36%  '
